$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 intermediate repeated-worker rows (old rows 18-21 for JESUS DAVID
# CACERES NIÑO periods 2506-2503, plus the old row 17 duplicate for period 2507);
# this shifts the former "last row" (old row 22, periods 2502) up to become the
# new row 17, and shifts the signature block (old rows 27-28) up to rows 22-23.
$ws.Rows("17:21").Delete()

# Update the Valor Mora total for the new 2-row table.
$ws.Range("E11").Value = 113880

# Cant. Trabajadores: 2 -> 1, Cant. Periodos: 6 -> 2
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 2

# Row 17 becomes the second statement line for MARTHA CECILIA VELEZ CALLE,
# period 2508 (new account-statement data, "parte 1").
$ws.Range("C17").Value = "1143374199"
$ws.Range("D17").Value = "MARTHA CECILIA VELEZ CALLE"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500
